# Ripple template update:
#  - add a new "Assay" worksheet at the end of the workbook containing the
#    default assay/dilution settings table
#  - (cosmetic) move the active selection on the "Patterns" sheet

$wb = $excel.ActiveWorkbook

# --- Add the new "Assay" sheet after the last existing sheet -----------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$assay = $wb.Worksheets.Add($null, $lastSheet)
$assay.Name = "Assay"

# --- Populate the Setting / Value table ---------------------------------
$settings = @(
    @("Setting", "Value"),
    @("DMSO Tolerance", 0.005),
    @("Well Volume (µL)", 25),
    @("Backfill (µL)", 10),
    @("Allowed Error", 0.1),
    @("Destination Replicates", 1),
    @("Use Intermediate Plates", 1),
    @("DMSO Normalization", 1)
)

for ($i = 0; $i -lt $settings.Length; $i++) {
    $row = $i + 1
    $assay.Cells.Item($row, 1).Value = $settings[$i][0]
    $assay.Cells.Item($row, 2).Value = $settings[$i][1]
}

$assay.Range("A1:B8").Select()

# --- Restore the Patterns sheet as the active tab / selection ----------
$patterns = $wb.Worksheets.Item("Patterns")
$patterns.Activate()
$patterns.Range("J24").Select()
